$d = $word.ActiveDocument

# Locate the last paragraph (the "Git restore ... this will restore the commit " one)
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)   # wdCollapseEnd -> collapse to end of paragraph (before the paragraph mark)
$r.MoveEnd(1, -1) # move before paragraph mark character, if needed

# Insert the extra run text right after the existing text, before the paragraph mark.
$endRange = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)
$endRange.InsertAfter("(do it before commiting )")

# Now add a new paragraph after this one for "Git reset ..." text
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$p2 = $lastPara.Range.InsertParagraphAfter()

$resetPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$resetRange = $resetPara.Range
$resetRange.Collapse(0)
$resetRange.InsertBefore("Git reset ")

$codeRange = $d.Range($resetPara.Range.End - 1, $resetPara.Range.End - 1)
$codeRange.InsertAfter("cf3cfe390d0c1a7c88f70ce3a892c667149c5d5f " + [char]8211 + " this is used when you want to uncommit the changes from a file and move them back to unstagged changes")
$codeRange.Font.Name = "Lucida Console"
$codeRange.Font.NameBi = "Lucida Console"
$codeRange.Font.Size = 9
$codeRange.Font.SizeBi = 9
$codeRange.Font.Color = 49087

# New paragraph: "Use git log after this to see the commites"
$resetPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$p3 = $resetPara.Range.InsertParagraphAfter()

$logPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$logRange = $logPara.Range
$logRange.Collapse(0)
$logRange.InsertAfter("Use git log after this to see the commites")
$logRange.Font.Name = "Lucida Console"
$logRange.Font.Size = 9
$logRange.Font.Color = 49087
